$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 9615686
$ws.Range("I2").Value = 142.71428
$ws.Range("J2").Value = 20833820
$ws.Range("K2").Value = 142.71428
$ws.Range("L2").Value = 20833820
$ws.Range("M2").Value = -29.71428
$ws.Range("N2").Value = -20834046

$ws.Range("H74").Value = 3894
$ws.Range("I74").Value = 3724.3635
$ws.Range("J74").Value = 4063.6365
$ws.Range("K74").Value = 3724.3635
$ws.Range("L74").Value = 4063.6365
$ws.Range("M74").Value = -2788.3635
$ws.Range("N74").Value = -5935.636500000001

$ws.Range("H77").Value = 3894
$ws.Range("I77").Value = 3724.3635
$ws.Range("J77").Value = 4063.6365
$ws.Range("K77").Value = 18621.8175
$ws.Range("L77").Value = 20318.1825
$ws.Range("M77").Value = -13941.8175
$ws.Range("N77").Value = -29678.1825

$ws.Range("H82").Value = 2163.2856
$ws.Range("I82").Value = 429.6
$ws.Range("J82").Value = 6497.5
$ws.Range("K82").Value = 1288.8
$ws.Range("L82").Value = 19492.5
$ws.Range("M82").Value = -882.8000000000002
$ws.Range("N82").Value = -20304.5

$ws.Range("H85").Value = 2163.2856
$ws.Range("I85").Value = 429.6
$ws.Range("J85").Value = 6497.5
$ws.Range("K85").Value = 1288.8
$ws.Range("L85").Value = 19492.5
$ws.Range("M85").Value = 115.1999999999998
$ws.Range("N85").Value = -22300.5

$ws.Range("H88").Value = 2338.4688
$ws.Range("J88").Value = 1234.875
$ws.Range("L88").Value = 1234.875
$ws.Range("N88").Value = -2046.875

$ws.Range("H91").Value = 2338.4688
$ws.Range("J91").Value = 1234.875
$ws.Range("L91").Value = 1234.875
$ws.Range("N91").Value = -4042.875

$ws.Range("H98").Value = 711.73914
$ws.Range("I98").Value = 711.73914
$ws.Range("K98").Value = 711.73914
$ws.Range("M98").Value = 786.26086

$ws.Range("H122").Value = 711.73914
$ws.Range("I122").Value = 711.73914
$ws.Range("K122").Value = 2135.21742
$ws.Range("M122").Value = 314.7825800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 11510
$ws.Range("I97").Value = 12013.333
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 12013.333
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -11517.333
$ws.Range("N97").Value = -10992

$ws.Range("H110").Value = 1997.5358
$ws.Range("I110").Value = 2250.4
$ws.Range("K110").Value = 2250.4
$ws.Range("M110").Value = -205.4000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1100.3334
$ws.Range("I22").Value = 301
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 301
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -128
$ws.Range("N22").Value = -1846

$ws.Range("H86").Value = 4764406.5
$ws.Range("I86").Value = 5716697
$ws.Range("J86").Value = 2955.2856
$ws.Range("K86").Value = 5716697
$ws.Range("L86").Value = 2955.2856
$ws.Range("M86").Value = -5715574
$ws.Range("N86").Value = -5201.2856

$ws.Range("H89").Value = 4764406.5
$ws.Range("I89").Value = 5716697
$ws.Range("J89").Value = 2955.2856
$ws.Range("K89").Value = 28583485
$ws.Range("L89").Value = 14776.428
$ws.Range("M89").Value = -28577869
$ws.Range("N89").Value = -26008.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 61485.555
$ws.Range("J4").Value = 6671.25
$ws.Range("L4").Value = 6671.25
$ws.Range("N4").Value = -6895.25

$ws.Range("H16").Value = 3312.5625
$ws.Range("I16").Value = 3888.889
$ws.Range("J16").Value = 2571.5715
$ws.Range("K16").Value = 3888.889
$ws.Range("L16").Value = 2571.5715
$ws.Range("M16").Value = -3601.889
$ws.Range("N16").Value = -3145.5715

$ws.Range("H22").Value = 235.44444
$ws.Range("I22").Value = 202.71428
$ws.Range("K22").Value = 202.71428
$ws.Range("M22").Value = 147.28572

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H58").Value = 3245.6365
$ws.Range("I58").Value = 750.2857
$ws.Range("J58").Value = 7612.5
$ws.Range("K58").Value = 750.2857
$ws.Range("L58").Value = 7612.5
$ws.Range("M58").Value = -547.2857
$ws.Range("N58").Value = -8018.5

$ws.Range("H86").Value = 337240.94
$ws.Range("I86").Value = 558190.4399999999
$ws.Range("J86").Value = 5816.6665
$ws.Range("K86").Value = 558190.4399999999
$ws.Range("L86").Value = 5816.6665
$ws.Range("M86").Value = -557067.4399999999
$ws.Range("N86").Value = -8062.6665

$ws.Range("H89").Value = 337240.94
$ws.Range("I89").Value = 558190.4399999999
$ws.Range("J89").Value = 5816.6665
$ws.Range("K89").Value = 2790952.2
$ws.Range("L89").Value = 29083.3325
$ws.Range("M89").Value = -2785336.2
$ws.Range("N89").Value = -40315.3325

$ws.Range("H105").Value = 713.625
$ws.Range("I105").Value = 672.7143
$ws.Range("K105").Value = 672.7143
$ws.Range("M105").Value = 1074.2857

$ws.Range("H113").Value = 3312.5625
$ws.Range("I113").Value = 3888.889
$ws.Range("J113").Value = 2571.5715
$ws.Range("K113").Value = 3888.889
$ws.Range("L113").Value = 2571.5715
$ws.Range("M113").Value = -1718.889
$ws.Range("N113").Value = -6911.5715

$ws.Range("H136").Value = 3245.6365
$ws.Range("I136").Value = 750.2857
$ws.Range("J136").Value = 7612.5
$ws.Range("K136").Value = 2250.8571
$ws.Range("L136").Value = 22837.5
$ws.Range("M136").Value = 299.1428999999998
$ws.Range("N136").Value = -27937.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 313.85715
$ws.Range("I8").Value = 313.85715
$ws.Range("K8").Value = 941.5714499999999
$ws.Range("M8").Value = -802.5714499999999

$ws.Range("H113").Value = 6439.8
$ws.Range("I113").Value = 533
$ws.Range("J113").Value = 15300
$ws.Range("K113").Value = 1599
$ws.Range("L113").Value = 45900
$ws.Range("M113").Value = 571
$ws.Range("N113").Value = -50240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5981
$ws.Range("I80").Value = 6826.25
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 6826.25
$ws.Range("L80").Value = 2600
$ws.Range("M80").Value = -5828.25
$ws.Range("N80").Value = -4596

$ws.Range("H83").Value = 5981
$ws.Range("I83").Value = 6826.25
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 34131.25
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = -29139.25
$ws.Range("N83").Value = -22984

$ws.Range("H113").Value = 1375.4706
$ws.Range("I113").Value = 821.5
$ws.Range("K113").Value = 821.5
$ws.Range("M113").Value = 1348.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2223166.2
$ws.Range("I22").Value = 3030983
$ws.Range("J22").Value = 1670
$ws.Range("K22").Value = 3030983
$ws.Range("L22").Value = 1670
$ws.Range("M22").Value = -3030688
$ws.Range("N22").Value = -2260

$ws.Range("H27").Value = 2223166.2
$ws.Range("I27").Value = 3030983
$ws.Range("J27").Value = 1670
$ws.Range("K27").Value = 3030983
$ws.Range("L27").Value = 1670
$ws.Range("M27").Value = -3030876
$ws.Range("N27").Value = -1884

$ws.Range("H82").Value = 7577406
$ws.Range("I82").Value = 15152530
$ws.Range("J82").Value = 2281.6667
$ws.Range("K82").Value = 15152530
$ws.Range("L82").Value = 2281.6667
$ws.Range("M82").Value = -15152169
$ws.Range("N82").Value = -3003.6667

$ws.Range("H85").Value = 7577406
$ws.Range("I85").Value = 15152530
$ws.Range("J85").Value = 2281.6667
$ws.Range("K85").Value = 15152530
$ws.Range("L85").Value = 2281.6667
$ws.Range("M85").Value = -15151282
$ws.Range("N85").Value = -4777.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9126.6875
$ws.Range("J2").Value = 9401.799999999999
$ws.Range("L2").Value = 9401.799999999999
$ws.Range("N2").Value = -9625.799999999999

$ws.Range("H62").Value = 255058.33
$ws.Range("I62").Value = 4987.5
$ws.Range("J62").Value = 755200
$ws.Range("K62").Value = 4987.5
$ws.Range("L62").Value = 755200
$ws.Range("M62").Value = -4363.5
$ws.Range("N62").Value = -756448

$ws.Range("H65").Value = 255058.33
$ws.Range("I65").Value = 4987.5
$ws.Range("J65").Value = 755200
$ws.Range("K65").Value = 24937.5
$ws.Range("L65").Value = 3776000
$ws.Range("M65").Value = -21817.5
$ws.Range("N65").Value = -3782240

$ws.Range("H100").Value = 40866.668
$ws.Range("I100").Value = 199533.33
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 399066.66
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -398525.66
$ws.Range("N100").Value = -3482

$ws.Range("H126").Value = 4254.8887
$ws.Range("I126").Value = 6015
$ws.Range("J126").Value = 734.6667
$ws.Range("K126").Value = 18045
$ws.Range("L126").Value = 2204.0001
$ws.Range("M126").Value = -15575
$ws.Range("N126").Value = -7144.0001
